$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 6).Value = 198
$ws1.Cells.Item(3, 6).Value = 3386
$ws1.Cells.Item(4, 6).Value = 249
$ws1.Cells.Item(5, 6).Value = 144
$ws1.Cells.Item(7, 6).Value = 1735
$ws1.Cells.Item(10, 6).Value = 380
$ws1.Cells.Item(12, 6).Value = 32
$ws1.Cells.Item(17, 6).Value = 33
$ws1.Cells.Item(20, 6).Value = 11
$ws1.Cells.Item(21, 6).Value = 29
$ws1.Cells.Item(24, 6).Value = 43
$ws1.Cells.Item(26, 6).Value = 403
$ws1.Cells.Item(27, 6).Value = 263
$ws1.Cells.Item(28, 6).Value = 118
$ws1.Cells.Item(29, 6).Value = 43
$ws1.Cells.Item(30, 6).Value = 19
$ws1.Cells.Item(32, 6).Value = 442
$ws1.Cells.Item(33, 6).Value = 2311
$ws1.Cells.Item(35, 6).Value = 53
$ws1.Cells.Item(36, 6).Value = 484
$ws1.Cells.Item(37, 6).Value = 569
$ws1.Cells.Item(40, 6).Value = 238
$ws1.Cells.Item(41, 6).Value = 358
$ws1.Cells.Item(42, 6).Value = 418
$ws1.Cells.Item(43, 6).Value = 542

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(2, 6).Value = 198
$ws4.Cells.Item(3, 6).Value = 3386
$ws4.Cells.Item(4, 6).Value = 249
$ws4.Cells.Item(5, 6).Value = 144
$ws4.Cells.Item(7, 6).Value = 1735
$ws4.Cells.Item(10, 6).Value = 380
$ws4.Cells.Item(12, 6).Value = 32
$ws4.Cells.Item(17, 6).Value = 33
$ws4.Cells.Item(20, 6).Value = 11
$ws4.Cells.Item(21, 6).Value = 29
$ws4.Cells.Item(24, 6).Value = 43
$ws4.Cells.Item(26, 6).Value = 403
$ws4.Cells.Item(27, 6).Value = 263
$ws4.Cells.Item(28, 6).Value = 118
$ws4.Cells.Item(29, 6).Value = 43
$ws4.Cells.Item(30, 6).Value = 19
$ws4.Cells.Item(32, 6).Value = 442
$ws4.Cells.Item(33, 6).Value = 2311
$ws4.Cells.Item(35, 6).Value = 53
$ws4.Cells.Item(36, 6).Value = 484
$ws4.Cells.Item(37, 6).Value = 569
$ws4.Cells.Item(40, 6).Value = 239
$ws4.Cells.Item(41, 6).Value = 358
$ws4.Cells.Item(42, 6).Value = 418
$ws4.Cells.Item(43, 6).Value = 542
